$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force text format for cells whose new numeric-looking values would
# otherwise be auto-converted to numbers by Excel's General format.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '27.234.65'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '1.571.89'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').Value = '211.33'
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('E6').Value = '  +0.67%  '
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('E8').Value = '  -0.31%  '
$ws.Range('D9').Value = '0.248'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('D12').Value = '1.795.59'
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('D13').Value = '1.582.08'
$ws.Range('E13').Value = '  +1.33%  '
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').Value = '27.183.14'
$ws.Range('E16').Value = '  +0.74%  '
$ws.Range('D17').Value = '62.22'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').Value = '216.25'
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('E22').Value = '  +1.40%  '
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').Value = '1.95'
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('D25').Value = '154.00'
$ws.Range('E25').Value = '  +0.19%  '
$ws.Range('D26').Value = '6.66'
$ws.Range('E26').Value = '  +0.72%  '
$ws.Range('D27').Value = '15.09'
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('E28').Value = '  +2.05%  '
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('E30').Value = '  +2.75%  '
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('E32').Value = '  +0.46%  '
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.450.02'
$ws.Range('E33').Value = '  +1.94%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '3.17'
$ws.Range('E34').Value = '  +1.63%  '
$ws.Range('D35').Value = '1.12'
$ws.Range('E35').Value = '  +7.36%  '
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('E37').Value = '  +1.15%  '
$ws.Range('E38').Value = '  +1.12%  '
$ws.Range('E39').Value = '  +0.65%  '
$ws.Range('E40').Value = '  +2.13%  '
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('D43').Value = '2.34'
$ws.Range('E43').Value = '  +0.53%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').Value = '64.50'
$ws.Range('E45').Value = '  -0.52%  '
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('D47').Value = '1.707.03'
$ws.Range('E47').Value = '  +0.69%  '
$ws.Range('D48').Value = '85.93'
$ws.Range('E48').Value = '  -1.65%  '
$ws.Range('E49').Value = '  +3.80%  '
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('D51').Value = '0.0959'
$ws.Range('E51').Value = '  +0.39%  '
